$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on price cells whose new values look numeric,
# so Excel stores them as text (matching the source data) instead of
# silently converting to floating point numbers.
$textCells = 'D5','D6','D10','D11','D12','D13','D14','D20','D21','D23','D24','D25','D27','D28','D29','D30','D31','D33','D34','D36','D37','D38','D39','D40','D42','D43','D44','D45','D47','D49','D51'
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '63.521.39'
$ws.Range('E2').Value = '  +2.29%  '

$ws.Range('D3').Value = '3.041.54'
$ws.Range('E3').Value = '  +0.95%  '

$ws.Range('E4').Value = '  +0.07%  '

$ws.Range('D5').Value = '553.37'
$ws.Range('E5').Value = '  +2.37%  '

$ws.Range('D6').Value = '141.03'
$ws.Range('E6').Value = '  +4.34%  '

$ws.Range('E7').Value = '  -0.13%  '

$ws.Range('D8').Value = '3.037.72'
$ws.Range('E8').Value = '  +1.02%  '

$ws.Range('E9').Value = '  +2.67%  '

$ws.Range('D10').Value = '0.155'
$ws.Range('E10').Value = '  +5.01%  '

$ws.Range('D11').Value = '6.05'
$ws.Range('E11').Value = '  -7.07%  '

$ws.Range('D12').Value = '0.473'
$ws.Range('E12').Value = '  +6.07%  '

$ws.Range('D13').Value = '0.0000229'
$ws.Range('E13').Value = '  +3.99%  '

$ws.Range('D14').Value = '34.54'
$ws.Range('E14').Value = '  +2.61%  '

$ws.Range('D15').Value = '3.534.64'
$ws.Range('E15').Value = '  +0.52%  '

$ws.Range('D16').Value = '63.564.74'
$ws.Range('E16').Value = '  +2.54%  '

$ws.Range('E17').Value = '  +1.68%  '

$ws.Range('D18').Value = '3.039.12'
$ws.Range('E18').Value = '  +0.91%  '

$ws.Range('E19').Value = '  +1.46%  '

$ws.Range('D20').Value = '476.39'
$ws.Range('E20').Value = '  +1.70%  '

$ws.Range('D21').Value = '13.93'
$ws.Range('E21').Value = '  +3.55%  '

$ws.Range('E22').Value = '  +2.74%  '

$ws.Range('D23').Value = '7.50'
$ws.Range('E23').Value = '  +5.86%  '

$ws.Range('D24').Value = '14.05'
$ws.Range('E24').Value = '  +12.96%  '

$ws.Range('D25').Value = '80.91'
$ws.Range('E25').Value = '  +2.06%  '

$ws.Range('E26').Value = '  -0.01%  '

$ws.Range('D27').Value = '2.77'
$ws.Range('E27').Value = '  +2.43%  '

$ws.Range('D28').Value = '7.88'
$ws.Range('E28').Value = '  +3.49%  '

$ws.Range('D29').Value = '2.02'
$ws.Range('E29').Value = '  +2.10%  '

$ws.Range('D30').Value = '1.00'
$ws.Range('E30').Value = '  +0.28%  '

$ws.Range('D31').Value = '26.12'
$ws.Range('E31').Value = '  +2.48%  '

$ws.Range('E32').Value = '  +0.40%  '

$ws.Range('D33').Value = '2.41'
$ws.Range('E33').Value = '  +2.97%  '

$ws.Range('D34').Value = '5.59'
$ws.Range('E34').Value = '  +1.06%  '

$ws.Range('E35').Value = '  +5.49%  '

$ws.Range('D36').Value = '54.62'
$ws.Range('E36').Value = '  +0.43%  '

$ws.Range('D37').Value = '0.0404'
$ws.Range('E37').Value = '  +2.85%  '

$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').Value = '436.12'
$ws.Range('E38').Value = '  -3.24%  '

$ws.Range('B39').Value = 'Hedera'
$ws.Range('C39').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D39').Value = '0.0804'
$ws.Range('E39').Value = '  +0.56%  '

$ws.Range('B40').Value = 'dogwifhat'
$ws.Range('C40').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D40').Value = '2.85'
$ws.Range('E40').Value = '  +17.71%  '

$ws.Range('D41').Value = '2.965.51'
$ws.Range('E41').Value = '  +0.57%  '

$ws.Range('D42').Value = '8.19'
$ws.Range('E42').Value = '  +2.61%  '

$ws.Range('D43').Value = '0.112'
$ws.Range('E43').Value = '  -2.82%  '

$ws.Range('D44').Value = '28.11'
$ws.Range('E44').Value = '  +5.54%  '

$ws.Range('D45').Value = '0.256'
$ws.Range('E45').Value = '  +3.15%  '

$ws.Range('E46').Value = '  -0.02%  '

$ws.Range('D47').Value = '2.10'
$ws.Range('E47').Value = '  +6.25%  '

$ws.Range('E48').Value = '  +3.17%  '

$ws.Range('D49').Value = '116.62'
$ws.Range('E49').Value = '  +1.77%  '

$ws.Range('D50').Value = '0.0₃0511'
$ws.Range('E50').Value = '  +4.74%  '

$ws.Range('D51').Value = '2.06'
$ws.Range('E51').Value = '  +2.92%  '

# Restore the default (Normal) style so these cells don't carry a
# leftover explicit number-format style compared to the original file.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}